# mindSpore/flora_quest.xlsx - "Add files via upload"
#
# Update the score for the second team member (B2: 30 -> 70), move the
# active cell selection to B2, and set the sheet's print/page setup
# (paper size = A4, orientation = portrait) which is how Excel ends up
# emitting a <pageSetup .../> element for the sheet on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2 changed from 30 to 70
$ws.Range("B2").Value = 70

# Selection moved from A3 to B2
$ws.Range("B2").Select()

# Page setup: paper size A4, portrait orientation
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
